$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns (and new Coin/Link cells) are stored as
# plain text, matching the original workbook where these are inline strings
# (e.g. "10.40" must stay "10.40", not become the number 10.4).
$ws.Range("B22:E51").NumberFormat = "@"
$ws.Range("D2:E21").NumberFormat = "@"

# Update rows 2-21: only Price (D) and Volume(1h) (E) columns change
$ws.Range("D2").Value = "60.596.74"
$ws.Range("E2").Value = "  +3.29%  "
$ws.Range("D3").Value = "2.660.75"
$ws.Range("E3").Value = "  +1.09%  "
$ws.Range("D4").Value = "0.996"
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").Value = "568.35"
$ws.Range("E5").Value = "  +5.91%  "
$ws.Range("D6").Value = "145.50"
$ws.Range("E6").Value = "  +1.58%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "0.606"
$ws.Range("E8").Value = "  +6.93%  "
$ws.Range("D9").Value = "6.84"
$ws.Range("E9").Value = "  -2.99%  "
$ws.Range("D10").Value = "0.105"
$ws.Range("E10").Value = "  +4.29%  "
$ws.Range("D11").Value = "0.143"
$ws.Range("E11").Value = "  +6.36%  "
$ws.Range("D12").Value = "0.341"
$ws.Range("E12").Value = "  +2.21%  "
$ws.Range("D13").Value = "3.104.38"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").Value = "60.217.35"
$ws.Range("E14").Value = "  +2.73%  "
$ws.Range("D15").Value = "21.73"
$ws.Range("E15").Value = "  +3.88%  "
$ws.Range("D16").Value = "2.646.83"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").Value = "0.0000136"
$ws.Range("E17").Value = "  +2.76%  "
$ws.Range("D18").Value = "4.57"
$ws.Range("E18").Value = "  +3.57%  "
$ws.Range("D19").Value = "343.49"
$ws.Range("E19").Value = "  +2.36%  "
$ws.Range("D20").Value = "10.40"
$ws.Range("E20").Value = "  +2.23%  "
$ws.Range("D21").Value = "6.35"
$ws.Range("E21").Value = "  +1.66%  "

# Rows 22-51: a new coin (LEO) is inserted at row 22, shifting all subsequent
# coins down by one row; the last row (Maker) is dropped off the bottom.
# Update Coin (B), Link (C), Price (D) and Volume(1h) (E) columns accordingly.
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").Value = "5.80"
$ws.Range("E22").Value = "  +0.46%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "66.54"
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").Value = "0.437"
$ws.Range("E25").Value = "  +5.17%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "0.166"
$ws.Range("E26").Value = "  +1.38%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "0.996"
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "7.31"
$ws.Range("E28").Value = "  +2.01%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0771"
$ws.Range("E29").Value = "  +4.26%  "
$ws.Range("B30").Value = "USDe"
$ws.Range("C30").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D30").Value = "0.997"
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "1.71"
$ws.Range("E31").Value = "  +3.79%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "6.13"
$ws.Range("E32").Value = "  +4.54%  "
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "156.23"
$ws.Range("E33").Value = "  +3.70%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "19.16"
$ws.Range("E34").Value = "  +2.13%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "4.09"
$ws.Range("E35").Value = "  +4.50%  "
$ws.Range("B36").Value = "SuiNetwork"
$ws.Range("C36").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D36").Value = "0.906"
$ws.Range("E36").Value = "  +6.99%  "
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").Value = "0.910"
$ws.Range("E37").Value = "  +11.73%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "1.16"
$ws.Range("E38").Value = "  +5.22%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "37.48"
$ws.Range("E39").Value = "  +0.89%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "1.50"
$ws.Range("E40").Value = "  +5.60%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").Value = "304.80"
$ws.Range("E41").Value = "  +8.01%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "3.66"
$ws.Range("E42").Value = "  +2.09%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "0.605"
$ws.Range("E44").Value = "  +0.54%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "0.0976"
$ws.Range("E45").Value = "  +4.34%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "0.0548"
$ws.Range("E46").Value = "  +3.03%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "19.36"
$ws.Range("E47").Value = "  +1.11%  "
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").Value = "10.64"
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "0.0235"
$ws.Range("E49").Value = "  +4.29%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "124.28"
$ws.Range("E50").Value = "  +11.03%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "4.70"
$ws.Range("E51").Value = "  +5.24%  "
